# Weekly refresh of "Fruta, Vega Monumental Concepción - Tuna" price records.
# Historical rows (2-14) are re-shuffled to their updated positions/values,
# and a new week's record (date 45084) is appended as row 15.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 44280
$ws.Range("N2").Value = 14000
$ws.Range("O2").Value = 15000
$ws.Range("P2").Value = 14500
$ws.Range("S2").Value = 806
$ws.Range("D3").Value = 44280
$ws.Range("L3").Value = 'Segunda'
$ws.Range("M3").Value = 50
$ws.Range("N3").Value = 12000
$ws.Range("O3").Value = 12000
$ws.Range("P3").Value = 12000
$ws.Range("Q3").Value = '$/caja 18 kilos'
$ws.Range("S3").Value = 667
$ws.Range("D4").Value = 45014
$ws.Range("M4").Value = 50
$ws.Range("N4").Value = 13000
$ws.Range("O4").Value = 14000
$ws.Range("P4").Value = 13600
$ws.Range("Q4").Value = '$/caja 18 kilos'
$ws.Range("S4").Value = 756
$ws.Range("L5").Value = 'Segunda'
$ws.Range("M5").Value = 20
$ws.Range("N5").Value = 10000
$ws.Range("O5").Value = 10000
$ws.Range("P5").Value = 10000
$ws.Range("S5").Value = 556
$ws.Range("D6").Value = 44699
$ws.Range("L6").Value = 'Primera'
$ws.Range("M6").Value = 100
$ws.Range("N6").Value = 20000
$ws.Range("O6").Value = 22000
$ws.Range("P6").Value = 21000
$ws.Range("S6").Value = 1167
$ws.Range("D7").Value = 44699
$ws.Range("L7").Value = 'Segunda'
$ws.Range("N7").Value = 18000
$ws.Range("O7").Value = 18000
$ws.Range("P7").Value = 18000
$ws.Range("S7").Value = 1000
$ws.Range("D8").Value = 44516
$ws.Range("N8").Value = 33000
$ws.Range("O8").Value = 34000
$ws.Range("P8").Value = 33500
$ws.Range("S8").Value = 1861
$ws.Range("D9").Value = 44687
$ws.Range("L9").Value = 'Primera'
$ws.Range("M9").Value = 100
$ws.Range("O9").Value = 19000
$ws.Range("P9").Value = 18500
$ws.Range("S9").Value = 1028
$ws.Range("D10").Value = 44819
$ws.Range("N10").Value = 25000
$ws.Range("O10").Value = 26000
$ws.Range("P10").Value = 25500
$ws.Range("Q10").Value = '$/caja 18 kilos granel'
$ws.Range("S10").Value = 1417
$ws.Range("D11").Value = 45084
$ws.Range("N11").Value = 20000
$ws.Range("O11").Value = 21000
$ws.Range("P11").Value = 20500
$ws.Range("Q11").Value = '$/caja 18 kilos granel'
$ws.Range("S11").Value = 1139
$ws.Range("D12").Value = 45044
$ws.Range("N12").Value = 17000
$ws.Range("O12").Value = 18000
$ws.Range("P12").Value = 17500
$ws.Range("S12").Value = 972
$ws.Range("D13").Value = 45002
$ws.Range("L13").Value = 'Primera'
$ws.Range("M13").Value = 100
$ws.Range("O13").Value = 13000
$ws.Range("P13").Value = 12500
$ws.Range("S13").Value = 694
$ws.Range("D14").Value = 44316
$ws.Range("M14").Value = 50
$ws.Range("N14").Value = 20000
$ws.Range("O14").Value = 20000
$ws.Range("P14").Value = 20000
$ws.Range("S14").Value = 1111
$ws.Range("A15").Value = 11
$ws.Range("B15").Value = 'Vega Monumental Concepción'
$ws.Range("C15").Value = 'Bíobío'
$ws.Range("D15").Value = 45030
$ws.Range("D15").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E15").Value = 8
$ws.Range("F15").Value = 'Fruta'
$ws.Range("G15").Value = 100107
$ws.Range("H15").Value = 'Otros'
$ws.Range("I15").Value = 100107011
$ws.Range("J15").Value = 'Tuna'
$ws.Range("K15").Value = 'Sin especificar'
$ws.Range("L15").Value = 'Primera'
$ws.Range("M15").Value = 100
$ws.Range("N15").Value = 15000
$ws.Range("O15").Value = 16000
$ws.Range("P15").Value = 15500
$ws.Range("Q15").Value = '$/caja 18 kilos granel'
$ws.Range("R15").Value = 'Provincia de Melipilla'
$ws.Range("S15").Value = 861
$ws.Range("T15").Value = 18
